# Pruebas sin ACK.xlsx - "Comprobacion de direccion y payload.Payload Variable"
# Adds a new worksheet "Retardos3" (with payload/ack timing data, mirroring the
# layout already used by "Retardos2" but with an extra "t3(us)" timing column),
# and updates the previously-active sheet's selection/ tab state accordingly.

$wb = $excel.ActiveWorkbook

# Helper: write a rectangular block of values starting at the given range
# (top-left anchored) using a real 2-D COM SAFEARRAY so Excel treats it as a
# single paste of values rather than one cell at a time.
function Set-Block {
    param($ws, [string]$range, $data)
    $rowCount = $data.Length
    $colCount = $data[0].Length
    $arr = New-Object 'object[,]' $rowCount, $colCount
    for ($r = 0; $r -lt $rowCount; $r++) {
        for ($c = 0; $c -lt $colCount; $c++) {
            $arr[$r, $c] = $data[$r][$c]
        }
    }
    $ws.Range($range).Value = $arr
}

# --- Select the range that will remain selected on Retardos2 once it stops
#     being the active tab, matching the target document state. ---
$ws2 = $wb.Worksheets.Item("Retardos2")
$ws2.Range("A5:N18").Select()

# --- Create the new worksheet right after "Retardos2" ---
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "Retardos3"

# ============================================================
# Row 5 - payload size headers (one per block)
# ============================================================
$ws3.Range("A5").Value = "Payload 12 bytes"
$ws3.Range("D5").Font.Bold = $true
$ws3.Range("E5").Value = "Payload 18 bytes"
$ws3.Range("I5").Value = "Payload 50 bytes"
$ws3.Range("M5").Value = "Payload 102 bytes"
$ws3.Range("Q5").Value = "Payload 114 bytes"
$ws3.Range("A5").Font.Bold = $true
$ws3.Range("E5").Font.Bold = $true
$ws3.Range("I5").Font.Bold = $true
$ws3.Range("M5").Font.Bold = $true
$ws3.Range("Q5").Font.Bold = $true

# ============================================================
# Row 7 - column headers (t1 / t2(uS) / t3(us)) for each block
# ============================================================
$ws3.Range("A7").Value = "t1"
$ws3.Range("B7").Value = "t2(uS)"
$ws3.Range("C7").Value = "t3(us)"
$ws3.Range("E7").Value = "t1"
$ws3.Range("F7").Value = "t2(uS)"
$ws3.Range("G7").Value = "t3(us)"
$ws3.Range("I7").Value = "t1"
$ws3.Range("J7").Value = "t2(uS)"
$ws3.Range("K7").Value = "t3(us)"
$ws3.Range("M7").Value = "t1"
$ws3.Range("N7").Value = "t2(uS)"
$ws3.Range("O7").Value = "t3(us)"
$ws3.Range("Q7").Value = "t1"
$ws3.Range("R7").Value = "t2(uS)"
$ws3.Range("S7").Value = "t3(us)"
$ws3.Range("A7:D7").Font.Bold = $true
$ws3.Range("E7:H7").Font.Bold = $true
$ws3.Range("I7:K7").Font.Bold = $true
$ws3.Range("M7:O7").Font.Bold = $true
$ws3.Range("Q7:S7").Font.Bold = $true

# ============================================================
# Rows 8-17 - raw sample data for each of the 5 blocks
# ============================================================
Set-Block $ws3 "A8:C17" @(
    @(0,2471,1673617),
    @(0,3429,1678018),
    @(0,2793,1675154),
    @(0,1836,1672268),
    @(0,2152,1670577),
    @(0,2478,1671548),
    @(0,2468,1671292),
    @(0,3430,1670623),
    @(0,1530,1664429),
    @(0,1520,1676304)
)

Set-Block $ws3 "E8:G17" @(
    @(0,2359,1656254),
    @(0,1734,1654697),
    @(0,3001,1654310),
    @(0,3312,1651881),
    @(0,1724,1656887),
    @(0,1743,1661782),
    @(0,3318,1653395),
    @(0,2673,1655381),
    @(0,3959,1652185),
    @(0,2351,1654923)
)

Set-Block $ws3 "I8:K17" @(
    @(0,3445,1652082),
    @(0,4733,1651979),
    @(0,4731,1664343),
    @(0,5046,1650478),
    @(0,4405,1652135),
    @(0,4400,1653861),
    @(0,3443,1651521),
    @(0,4092,1657727),
    @(0,4084,1652406),
    @(0,3122,1652254)
)

Set-Block $ws3 "M8:O17" @(
    @(0,4899,1641943),
    @(0,5217,1652252),
    @(0,6811,1652011),
    @(0,6496,1650900),
    @(0,6179,1652319),
    @(0,6816,1648933),
    @(0,4581,1651365),
    @(0,6175,1651273),
    @(0,5859,1650342),
    @(0,5848,1650682)
)

Set-Block $ws3 "Q8:S17" @(
    @(0,2489,1),
    @(0,1850,3),
    @(0,2490,4),
    @(0,3448,123),
    @(0,2513,3),
    @(0,3129,1323),
    @(0,3774,4),
    @(0,1547,3),
    @(0,2485,2),
    @(0,2489,1)
)

# ============================================================
# Row 18 - "Promedio" label + AVERAGE() formulas per block
# ============================================================
$ws3.Range("A18").Value = "Promedio"
$ws3.Range("E18").Value = "Promedio"
$ws3.Range("I18").Value = "Promedio"
$ws3.Range("M18").Value = "Promedio"
$ws3.Range("Q18").Value = "Promedio"
$ws3.Range("A18").Font.Bold = $true
$ws3.Range("D18").Font.Bold = $true
$ws3.Range("E18").Font.Bold = $true
$ws3.Range("I18").Font.Bold = $true
$ws3.Range("M18").Font.Bold = $true
$ws3.Range("Q18").Font.Bold = $true

$ws3.Range("B18").Formula = "=AVERAGE(B8:B17)"
$ws3.Range("C18").Formula = "=AVERAGE(C8:C17)"
$ws3.Range("F18").Formula = "=AVERAGE(F8:F17)"
$ws3.Range("G18").Formula = "=AVERAGE(G8:G17)"
$ws3.Range("J18").Formula = "=AVERAGE(J8:J17)"
$ws3.Range("K18").Formula = "=AVERAGE(K8:K17)"
$ws3.Range("N18").Formula = "=AVERAGE(N8:N17)"
$ws3.Range("O18").Formula = "=AVERAGE(O8:O17)"
$ws3.Range("R18").Formula = "=AVERAGE(R8:R17)"
$ws3.Range("S18").Formula = "=AVERAGE(S8:S17)"

# ============================================================
# View state: make Retardos3 the active/visible tab, scrolled down a
# little with P17 selected, matching the saved workbook view.
# ============================================================
$ws3.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
$ws3.Range("P17").Select()
